# Daily attendance processing - 2025-12-22 22:57:20
#
# The "Recorded By" column (G) lists the users/services that recorded each
# attendance session, as a comma-separated string (e.g. "System, someone@x.com").
# This pass re-normalizes the ordering of that list for every data row on the
# "Session Analysis Results" sheet, flipping it end-to-end so the last-listed
# recorder comes first and the first-listed recorder comes last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if (($val -ne $null) -and ($val -is [string]) -and ($val.Contains(", "))) {
        $parts = $val -split ", "
        $reversed = $parts[($parts.Count - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
    }
}
